$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1,1).Value = 'Datos actualizados a 26 de Mayo de 2020 a las 01:05'

# Row 4
$ws.Cells.Item(4,2).Value = 1705668
$ws.Cells.Item(4,3).Value = 19232
$ws.Cells.Item(4,4).Value = 462640
$ws.Cells.Item(4,5).Value = 1143246
$ws.Cells.Item(4,7).Value = 482
$ws.Cells.Item(4,8).Value = 99782

# Row 5
$ws.Cells.Item(5,2).Value = 374898
$ws.Cells.Item(5,3).Value = 11280
$ws.Cells.Item(5,4).Value = 153833
$ws.Cells.Item(5,5).Value = 197592
$ws.Cells.Item(5,7).Value = 757
$ws.Cells.Item(5,8).Value = 23473

# Row 15
$ws.Cells.Item(15,2).Value = 123979
$ws.Cells.Item(15,3).Value = 4020
$ws.Cells.Item(15,4).Value = 50949
$ws.Cells.Item(15,5).Value = 69401
$ws.Cells.Item(15,7).Value = 173
$ws.Cells.Item(15,8).Value = 3629

# Row 16
$ws.Cells.Item(16,2).Value = 85711
$ws.Cells.Item(16,3).Value = 1012
$ws.Cells.Item(16,4).Value = 44638
$ws.Cells.Item(16,5).Value = 34528
$ws.Cells.Item(16,7).Value = 121
$ws.Cells.Item(16,8).Value = 6545

# Row 43
$ws.Cells.Item(43,2).Value = 16581
$ws.Cells.Item(43,3).Value = 31
$ws.Cells.Item(43,4).Value = 13612
$ws.Cells.Item(43,5).Value = 2139
$ws.Cells.Item(43,7).Value = 10
$ws.Cells.Item(43,8).Value = 830

# Row 51
$ws.Cells.Item(51,1).Value = 'Panama'
$ws.Cells.Item(51,2).Value = 11183
$ws.Cells.Item(51,3).Value = 257
$ws.Cells.Item(51,4).Value = 6279
$ws.Cells.Item(51,5).Value = 4594
$ws.Cells.Item(51,7).Value = 4
$ws.Cells.Item(51,8).Value = 310

# Row 52
$ws.Cells.Item(52,1).Value = 'Afganistan'
$ws.Cells.Item(52,2).Value = 11173
$ws.Cells.Item(52,3).Value = 591
$ws.Cells.Item(52,4).Value = 1097
$ws.Cells.Item(52,5).Value = 9857
$ws.Cells.Item(52,7).Value = 1
$ws.Cells.Item(52,8).Value = 219

# Row 54
$ws.Cells.Item(54,2).Value = 9002
$ws.Cells.Item(54,3).Value = 47
$ws.Cells.Item(54,4).Value = 6182
$ws.Cells.Item(54,5).Value = 2503

# Row 57
$ws.Cells.Item(57,2).Value = 8364
$ws.Cells.Item(57,3).Value = 12
$ws.Cells.Item(57,5).Value = 402

# Row 58
$ws.Cells.Item(58,2).Value = 8068
$ws.Cells.Item(58,3).Value = 229
$ws.Cells.Item(58,4).Value = 2311
$ws.Cells.Item(58,5).Value = 5524
$ws.Cells.Item(58,7).Value = 7
$ws.Cells.Item(58,8).Value = 233

# Row 120
$ws.Cells.Item(120,1).Value = 'Sudan del Sur'
$ws.Cells.Item(120,2).Value = 806
$ws.Cells.Item(120,3).Value = 151
$ws.Cells.Item(120,4).Value = 6
$ws.Cells.Item(120,5).Value = 792
$ws.Cells.Item(120,8).Value = 8

# Row 121
$ws.Cells.Item(121,1).Value = 'Uruguay'
$ws.Cells.Item(121,2).Value = 787
$ws.Cells.Item(121,3).Value = 18
$ws.Cells.Item(121,4).Value = 629
$ws.Cells.Item(121,5).Value = 136
$ws.Cells.Item(121,8).Value = 22

# Row 122
$ws.Cells.Item(122,1).Value = 'Principado de Andorra'
$ws.Cells.Item(122,2).Value = 763
$ws.Cells.Item(122,3).Value = 1
$ws.Cells.Item(122,4).Value = 663
$ws.Cells.Item(122,5).Value = 49
$ws.Cells.Item(122,7).Value = 0
$ws.Cells.Item(122,8).Value = 51

# Row 123
$ws.Cells.Item(123,1).Value = 'Sierra Leona'
$ws.Cells.Item(123,2).Value = 735
$ws.Cells.Item(123,3).Value = 28
$ws.Cells.Item(123,4).Value = 293
$ws.Cells.Item(123,5).Value = 400
$ws.Cells.Item(123,7).Value = 2
$ws.Cells.Item(123,8).Value = 42

# Row 124
$ws.Cells.Item(124,1).Value = 'Georgia'
$ws.Cells.Item(124,2).Value = 731
$ws.Cells.Item(124,3).Value = 1
$ws.Cells.Item(124,4).Value = 526
$ws.Cells.Item(124,5).Value = 193
$ws.Cells.Item(124,8).Value = 12

# Row 125
$ws.Cells.Item(125,1).Value = 'Crucero'
$ws.Cells.Item(125,2).Value = 712
$ws.Cells.Item(125,3).Value = 0
$ws.Cells.Item(125,4).Value = 651
$ws.Cells.Item(125,5).Value = 48
$ws.Cells.Item(125,8).Value = 13

# Row 126
$ws.Cells.Item(126,1).Value = 'Jordania'
$ws.Cells.Item(126,2).Value = 711
$ws.Cells.Item(126,3).Value = 3
$ws.Cells.Item(126,4).Value = 479
$ws.Cells.Item(126,5).Value = 223
$ws.Cells.Item(126,7).Value = 0
$ws.Cells.Item(126,8).Value = 9

# Row 127
$ws.Cells.Item(127,1).Value = 'Republica del Chad'
$ws.Cells.Item(127,2).Value = 687
$ws.Cells.Item(127,3).Value = 12
$ws.Cells.Item(127,4).Value = 244
$ws.Cells.Item(127,5).Value = 382
$ws.Cells.Item(127,8).Value = 61

# Row 128
$ws.Cells.Item(128,1).Value = 'Nepal'
$ws.Cells.Item(128,2).Value = 682
$ws.Cells.Item(128,3).Value = 79
$ws.Cells.Item(128,4).Value = 112
$ws.Cells.Item(128,5).Value = 566
$ws.Cells.Item(128,7).Value = 1
$ws.Cells.Item(128,8).Value = 4

# Row 129
$ws.Cells.Item(129,1).Value = 'San Marino'
$ws.Cells.Item(129,2).Value = 666
$ws.Cells.Item(129,3).Value = 1
$ws.Cells.Item(129,4).Value = 270
$ws.Cells.Item(129,5).Value = 354
$ws.Cells.Item(129,8).Value = 42

# Row 197
$ws.Cells.Item(197,1).Value = 'Curazao'
$ws.Cells.Item(197,3).Value = 1
$ws.Cells.Item(197,4).Value = 14
$ws.Cells.Item(197,8).Value = 1

# Row 198
$ws.Cells.Item(198,1).Value = 'Fiyi'
$ws.Cells.Item(198,3).Value = 0
$ws.Cells.Item(198,4).Value = 15
$ws.Cells.Item(198,8).Value = 0

# Row 199
$ws.Cells.Item(199,1).Value = 'Nueva Caledonia'

# Row 201
$ws.Cells.Item(201,1).Value = 'Santa Lucia'

# Row 215
$ws.Cells.Item(215,1).Value = 'Bonaire, San Eustaquio y Saba'

# Row 216
$ws.Cells.Item(216,1).Value = 'San Bartolome'
